$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update/rename activities whose label changed but keep their dates ---
# (column letters below are the ORIGINAL layout - B:F - the whole block is
#  shifted left into A:E at the very end of the script)
$ws.Range("B6").Value = "Betrachtung der Daten "
$ws.Range("B5").Value = "Fragestellung überlegen "
$ws.Range("B7").Value = "Skript zum Einlesen der Daten + Vorbereitung/ weitere Recherche"

# --- 2. Update the "Version" date stamp (row 2) ---
$ws.Range("F2").Value = 44342

# --- 3. Update duration/end date of the last existing activity (row 11) ---
$ws.Range("D11").Value = 84
$ws.Range("E11").Value = 44645

# --- 4. Add the new final activity (row 12), reusing the existing date style ---
$ws.Range("C4").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("B12").Value = "Abgabe"
$ws.Range("C12").Value = 44646
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 44650

# --- 5. Shift the whole table one column to the left (B:F -> A:E) ---
$ws.Columns.Item(1).Delete()

# --- 6. Resize the (now) text column A and the date column B ---
$ws.Columns.Item(1).ColumnWidth = 53.6
$ws.Columns.Item(2).ColumnWidth = 10.1

# --- 7. Switch the page to landscape orientation ---
$ws.PageSetup.Orientation = 2

# --- 8. Restore the selection to the new last cell ---
$ws.Range("D12").Select()
